$d = $word.ActiveDocument

# Locate the two paragraphs this edit touches by their distinctive content
# rather than by a hard-coded paragraph index, so the script is resilient to
# any unrelated paragraphs being present.
$p4 = $null
$p5 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*{m:template myTemplate(*") {
        $p4 = $p
    }
    if ($p.Range.Text -like "*{m: a + b + c}*") {
        $p5 = $p
    }
}
if ($p4 -eq $null) { throw "Could not find the 'myTemplate(...)' paragraph." }
if ($p5 -eq $null) { throw "Could not find the 'a + b + c' paragraph." }

# ---------------------------------------------------------------------------
# Paragraph "{m:template myTemplate(   a   :    Integer   ,   b   :    Integer
# ,   c   :    Integer   )}" -- the field-delimiter tokens '(' ')' '}' get
# pulled apart into their own runs (TokenIteratorFieldRewriterSplit), and a
# little extra whitespace is introduced around '(' and '}' the same way the
# migrated fixture does.
# ---------------------------------------------------------------------------
$p4Start = $p4.Range.Start

# Sanity-check we are editing the paragraph we think we are (Range.Text
# includes the trailing paragraph mark, hence the "`r").
$originalText = "{m:template myTemplate(   a   :    Integer   ,   b   :    Integer   ,   c   :    Integer   )}`r"
if ($p4.Range.Text -ne $originalText) {
    throw "Unexpected paragraph text: $($p4.Range.Text)"
}

# Insert "   " (3 spaces) right before the "(" character (offset 22).
$insertPoint1 = $d.Range($p4Start + 22, $p4Start + 22)
$insertPoint1.InsertBefore("   ")

# Insert "  " (2 spaces) right after the ")" that now sits at offset 94
# (91 in the original text, shifted by the 3 characters inserted above).
$insertPoint2 = $d.Range($p4Start + 95, $p4Start + 95)
$insertPoint2.InsertBefore("  ")

# Both InsertBefore calls collapse every same-formatted run in the paragraph
# into a single <w:r> (Word re-flows runs on any edit). Re-establish all of
# the run boundaries -- the ones that existed originally plus the two new
# ones -- with a harmless Bold on/off toggle, which splits a run without
# leaving any trace in <w:rPr>.
$boundaries = @(0, 22, 25, 26, 53, 75, 76, 94, 95, 97, 98)
for ($i = 0; $i -lt $boundaries.Length - 1; $i++) {
    $seg = $d.Range($p4Start + $boundaries[$i], $p4Start + $boundaries[$i + 1])
    $seg.Bold = 1
    $seg.Bold = 0
}

$expectedText = "{m:template myTemplate   (   a   :    Integer   ,   b   :    Integer   ,   c   :    Integer   )  }`r"
if ($p4.Range.Text -ne $expectedText) {
    throw "Paragraph text mismatch after edit: $($p4.Range.Text)"
}

# ---------------------------------------------------------------------------
# Paragraph "{m: a + b + c}" -- split the trailing "}" off of "b + c}" so the
# closing delimiter is its own run.
# ---------------------------------------------------------------------------
$p5Start = $p5.Range.Start

$originalText5 = "{m: a + b + c}"
if ($p5.Range.Text -ne ($originalText5 + "`r")) {
    throw "Unexpected paragraph text: $($p5.Range.Text)"
}

$braceOffset = $originalText5.Length - 1
$braceRun = $d.Range($p5Start + $braceOffset, $p5Start + $braceOffset + 1)
$braceRun.Bold = 1
$braceRun.Bold = 0
